$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 7 with the new component data, copying formatting from row 6
$ws.Range("A6:F6").Copy() | Out-Null
$ws.Range("A7:F7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A7").Value = "L34A196"
$ws.Range("B7").Value = "2002-1492"
$ws.Range("C7").Value = "Standard"
$ws.Range("D7").Value = "N/A"
$ws.Range("E7").Value = $null
$ws.Range("F7").Value = "SGD"

$ws.Range("D13").Select() | Out-Null
